# Update curso.xlsx: change the short "course code" values in column A
# for the ESO rows (1-12) from the "ESOxY" form (e.g. ESO1A) to the
# "xESOY" form (e.g. 1ESOA), and move the active selection to A12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codes = @{
    1  = "1ESOA"
    2  = "1ESOB"
    3  = "1ESOC"
    4  = "2ESOA"
    5  = "2ESOB"
    6  = "2ESOC"
    7  = "3ESOA"
    8  = "3ESOB"
    9  = "3ESOC"
    10 = "4ESOA"
    11 = "4ESOB"
    12 = "4ESOC"
}

foreach ($row in $codes.Keys) {
    $ws.Cells.Item($row, 1).Value = $codes[$row]
}

# Move the selection to A12, matching the saved sheetView state.
$ws.Range("A12").Select()
